$wb = $excel.ActiveWorkbook

# Sheet ALC row 18 (@@ -1520,22 +1520,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 420.5
$ws.Range("I18").Value = 420.5
$ws.Range("K18").Value = 420.5
$ws.Range("M18").Value = -136.5

# Sheet ALC row 38 (@@ -2485,22 +2485,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 977.4286
$ws.Range("I38").Value = 62
$ws.Range("K38").Value = 186
$ws.Range("M38").Value = 186

# Sheet ALC row 51 (@@ -3137,25 +3137,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9976.571
$ws.Range("I51").Value = 12965.25
$ws.Range("J51").Value = 5991.6665
$ws.Range("K51").Value = 12965.25
$ws.Range("L51").Value = 5991.6665
$ws.Range("M51").Value = -12481.25
$ws.Range("N51").Value = -6959.6665

# Sheet ALC row 58 (@@ -3489,22 +3489,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 584
$ws.Range("I58").Value = 618.9091
$ws.Range("J58").Value = 200
$ws.Range("K58").Value = 1856.7273
$ws.Range("L58").Value = 600
$ws.Range("M58").Value = -1706.7273
$ws.Range("N58").Value = -900

# Sheet ALC row 64 (@@ -3786,25 +3789,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5047.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 5047.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 5047.5
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -5543.5

# Sheet ALC row 67 (@@ -3939,25 +3939,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5047.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 5047.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 5047.5
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -6763.5

# Sheet ALC row 88 (@@ -4995,22 +4992,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1956.2273
$ws.Range("I88").Value = 2134.3333
$ws.Range("K88").Value = 2134.3333
$ws.Range("M88").Value = -1728.3333

# Sheet ALC row 91 (@@ -5145,22 +5142,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1956.2273
$ws.Range("I91").Value = 2134.3333
$ws.Range("K91").Value = 2134.3333
$ws.Range("M91").Value = -730.3332999999998

# Sheet ALC row 111 (@@ -6155,25 +6152,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1084.8
$ws.Range("I111").Value = 1084.8
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3254.4
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -187.3999999999996
$ws.Range("N111").ClearContents()

# Sheet ALC row 112 (@@ -6207,25 +6201,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1992.0952
$ws.Range("J112").Value = 1988
$ws.Range("L112").Value = 5964
$ws.Range("N112").Value = -8180

# Sheet ALC row 113 (@@ -6259,22 +6253,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4424.1665
$ws.Range("I113").Value = 3382.8333
$ws.Range("K113").Value = 3382.8333
$ws.Range("M113").Value = -128.8332999999998

# Sheet ALC row 141 (@@ -7640,25 +7634,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3649.087
$ws.Range("J141").Value = 1001.6667
$ws.Range("L141").Value = 3005.0001
$ws.Range("N141").Value = -13365.0001

# Sheet ARM row 2 (@@ -7786,25 +7780,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 925.2857
$ws.Range("I2").Value = 895.8
$ws.Range("J2").Value = 999
$ws.Range("K2").Value = 895.8
$ws.Range("L2").Value = 999
$ws.Range("M2").Value = -782.8
$ws.Range("N2").Value = -1225

# Sheet ARM row 88 (@@ -11946,25 +11940,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2048
$ws.Range("I88").Value = 1481.5
$ws.Range("J88").Value = 2173.889
$ws.Range("K88").Value = 1481.5
$ws.Range("L88").Value = 2173.889
$ws.Range("M88").Value = -1075.5
$ws.Range("N88").Value = -2985.889

# Sheet ARM row 91 (@@ -12093,25 +12087,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2048
$ws.Range("I91").Value = 1481.5
$ws.Range("J91").Value = 2173.889
$ws.Range("K91").Value = 1481.5
$ws.Range("L91").Value = 2173.889
$ws.Range("M91").Value = -77.5
$ws.Range("N91").Value = -4981.889

# Sheet ARM row 110 (@@ -13027,22 +13021,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1791.16
$ws.Range("I110").Value = 1751.4286
$ws.Range("K110").Value = 1751.4286
$ws.Range("M110").Value = 293.5714

# Sheet ARM row 116 (@@ -13318,25 +13312,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 925.2857
$ws.Range("I116").Value = 895.8
$ws.Range("J116").Value = 999
$ws.Range("K116").Value = 895.8
$ws.Range("L116").Value = 999
$ws.Range("M116").Value = 1398.2
$ws.Range("N116").Value = -5587

# Sheet BSM row 3 (@@ -14714,25 +14708,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 925.2857
$ws.Range("I3").Value = 895.8
$ws.Range("J3").Value = 999
$ws.Range("K3").Value = 895.8
$ws.Range("L3").Value = 999
$ws.Range("M3").Value = -781.8
$ws.Range("N3").Value = -1227

# Sheet BSM row 20 (@@ -15547,25 +15541,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2453.2144
$ws.Range("I20").Value = 2162
$ws.Range("J20").Value = 2977.4
$ws.Range("K20").Value = 2162
$ws.Range("L20").Value = 2977.4
$ws.Range("M20").Value = -1915
$ws.Range("N20").Value = -3471.4

# Sheet BSM row 80 (@@ -18460,25 +18454,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 350.25
$ws.Range("I80").Value = 313.5
$ws.Range("J80").Value = 387
$ws.Range("K80").Value = 313.5
$ws.Range("L80").Value = 387
$ws.Range("M80").Value = 684.5
$ws.Range("N80").Value = -2383

# Sheet BSM row 83 (@@ -18613,25 +18607,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 350.25
$ws.Range("I83").Value = 313.5
$ws.Range("J83").Value = 387
$ws.Range("K83").Value = 1567.5
$ws.Range("L83").Value = 1935
$ws.Range("M83").Value = 3424.5
$ws.Range("N83").Value = -11919

# Sheet BSM row 99 (@@ -19394,22 +19388,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 75232.07000000001
$ws.Range("I99").Value = 103324.9
$ws.Range("K99").Value = 103324.9
$ws.Range("M99").Value = -101826.9

# Sheet BSM row 107 (@@ -19789,25 +19783,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1779.2858
$ws.Range("J107").Value = 3998.5
$ws.Range("L107").Value = 3998.5
$ws.Range("N107").Value = -7838.5

# Sheet CRP row 7 (@@ -21804,25 +21798,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 795.7857
$ws.Range("J7").Value = 619.75
$ws.Range("L7").Value = 619.75
$ws.Range("N7").Value = -845.75

# Sheet CRP row 22 (@@ -22524,25 +22518,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1600.7
$ws.Range("I22").Value = 296.8
$ws.Range("J22").Value = 2904.6
$ws.Range("K22").Value = 296.8
$ws.Range("L22").Value = 2904.6
$ws.Range("M22").Value = 53.19999999999999
$ws.Range("N22").Value = -3604.6

# Sheet CRP row 31 (@@ -22962,25 +22956,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3789.6924
$ws.Range("I31").Value = 3137.25
$ws.Range("J31").Value = 4833.6
$ws.Range("K31").Value = 3137.25
$ws.Range("L31").Value = 4833.6
$ws.Range("M31").Value = -2842.25
$ws.Range("N31").Value = -5423.6

# Sheet CRP row 34 (@@ -23109,25 +23103,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3789.6924
$ws.Range("I34").Value = 3137.25
$ws.Range("J34").Value = 4833.6
$ws.Range("K34").Value = 3137.25
$ws.Range("L34").Value = 4833.6
$ws.Range("M34").Value = -2935.25
$ws.Range("N34").Value = -5237.6

# Sheet CRP row 62 (@@ -24481,23 +24475,20 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4995
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Sheet CRP row 65 (@@ -24625,23 +24616,20 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4995
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# Sheet CRP row 134 (@@ -27958,22 +27946,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 41691.23
$ws.Range("I134").Value = 47826.047
$ws.Range("K134").Value = 143478.141
$ws.Range("M134").Value = -140943.141

# Sheet CUL row 13 (@@ -29001,22 +28989,19 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

# Sheet CUL row 34 (@@ -30036,25 +30021,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1922.3334
$ws.Range("J34").Value = 1994.6666
$ws.Range("L34").Value = 5983.9998
$ws.Range("N34").Value = -6151.9998

# Sheet CUL row 39 (@@ -30287,22 +30272,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 9196.6
$ws.Range("J39").Value = 9196.6
$ws.Range("L39").Value = 27589.8
$ws.Range("N39").Value = -28177.8

# Sheet CUL row 55 (@@ -31083,25 +31068,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4398.5
$ws.Range("I55").Value = 4999.5
$ws.Range("J55").Value = 3797.5
$ws.Range("K55").Value = 14998.5
$ws.Range("L55").Value = 11392.5
$ws.Range("M55").Value = -14821.5
$ws.Range("N55").Value = -11746.5

# Sheet CUL row 128 (@@ -34723,22 +34708,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 153997.5
$ws.Range("I128").Value = 153997.5
$ws.Range("K128").Value = 461992.5
$ws.Range("M128").Value = -457012.5

# Sheet CUL row 140 (@@ -35335,22 +35320,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1965.0714
$ws.Range("I140").Value = 1965.0714
$ws.Range("K140").Value = 5895.2142
$ws.Range("M140").Value = -715.2142000000003

# Sheet GSM row 70 (@@ -38802,22 +38787,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19667
$ws.Range("I70").Value = 19667
$ws.Range("K70").Value = 19667
$ws.Range("M70").Value = -19397

# Sheet GSM row 73 (@@ -38943,22 +38928,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 19667
$ws.Range("I73").Value = 19667
$ws.Range("K73").Value = 19667
$ws.Range("M73").Value = -18731

# Sheet GSM row 126 (@@ -41522,22 +41507,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 9604.4
$ws.Range("I126").Value = 7816.7
$ws.Range("K126").Value = 23450.1
$ws.Range("M126").Value = -20980.1

# Sheet LTW row 46 (@@ -44517,25 +44502,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10510.375
$ws.Range("I46").Value = 24694.75
$ws.Range("J46").Value = 5782.25
$ws.Range("K46").Value = 24694.75
$ws.Range("L46").Value = 5782.25
$ws.Range("M46").Value = -24506.75
$ws.Range("N46").Value = -6158.25

# Sheet LTW row 55 (@@ -44949,25 +44934,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1624.6364
$ws.Range("I55").Value = 1513.3334
$ws.Range("J55").Value = 1758.2
$ws.Range("K55").Value = 1513.3334
$ws.Range("L55").Value = 1758.2
$ws.Range("M55").Value = -1340.3334
$ws.Range("N55").Value = -2104.2

# Sheet LTW row 107 (@@ -47503,22 +47488,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 2500
$ws.Range("I107").Value = 2500
$ws.Range("K107").Value = 2500
$ws.Range("M107").Value = -580

# Sheet WVR row 62 (@@ -52183,25 +52168,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 75490.17999999999
$ws.Range("I62").Value = 8610.5
$ws.Range("J62").Value = 134938.78
$ws.Range("K62").Value = 8610.5
$ws.Range("L62").Value = 134938.78
$ws.Range("M62").Value = -7986.5
$ws.Range("N62").Value = -136186.78

# Sheet WVR row 65 (@@ -52327,25 +52312,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 75490.17999999999
$ws.Range("I65").Value = 8610.5
$ws.Range("J65").Value = 134938.78
$ws.Range("K65").Value = 43052.5
$ws.Range("L65").Value = 674693.9
$ws.Range("M65").Value = -39932.5
$ws.Range("N65").Value = -680933.9

# Sheet WVR row 132 (@@ -55580,22 +55565,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 45837.316
$ws.Range("I132").Value = 66655.37
$ws.Range("K132").Value = 199966.11
$ws.Range("M132").Value = -197436.11
